$p = $ppt.ActivePresentation
Write-Host "--- Presentation members ---"
Write-Host (Get-Member -InputObject $p)
